$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$old1 = [char]0x2705 + " 1000 Bs = 8.84 = 36862.48 pesos"
$new1 = [char]0x2705 + " 1000 Bs = 8.85 = 36652.5 pesos"
$old2 = [char]0x2705 + " 36862.48 pesos = 8.84 = 957.28 Bs"
$new2 = [char]0x2705 + " 36652.5 pesos = 8.78 = 939.4 Bs"

$text = $wsHoja1.Range("A1").Value2
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update N10/O10/N12/O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 112.98
$wsTasas.Range("O10").Value = 4141
$wsTasas.Range("N12").Value = 4174.8
$wsTasas.Range("O12").Value = 107
